$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Then_Goto" / "Else_Goto" headers to "Then_Question" / "Else_Question"
$ws.Range("I1").Value = "Then_Question"
$ws.Range("J1").Value = "Else_Question"

# Update the active selection to J1 to match the author's saved state
$ws.Range("J1").Select()
